$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 57 (2025-05): retained_customers and retention_rate changed
$ws.Cells.Item(57, 2).Value = 143
$ws.Cells.Item(57, 4).Value = 68.75

# Add new row 58 (2025-06)
$ws.Cells.Item(58, 1).Value = "2025-06"
$ws.Cells.Item(58, 2).Value = 2
$ws.Cells.Item(58, 3).Value = 216
$ws.Cells.Item(58, 4).Value = 0.9259259259259258
